$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: single value correction ---
$ws.Range("B9").Value = 91833

# --- Rows 15 <-> 16: swap the two "Tretåig hackspett" records ---
# Row 15 becomes what row 16 used to be
$ws.Range("A15").Value = 131314795
$ws.Range("Q15").Value = 497547
$ws.Range("R15").Value = 6980312
$ws.Range("AC15").Value = "Ringhack, färska, på gran."
$ws.Range("AM15").Value = ""
$ws.Range("AO15").Value = "Picea abies"

# Row 16 becomes what row 15 used to be
$ws.Range("A16").Value = 131314696
$ws.Range("Q16").Value = 497570
$ws.Range("R16").Value = 6980217
$ws.Range("AC16").Value = "Ringhack, färska, på en gran."
$ws.Range("AM16").Value = "Trädstam på levande träd"
$ws.Range("AO16").Value = "Stem on living tree # Picea abies"

# --- Rows 20 <-> 22: swap the "Garnlav" and "Tretåig hackspett" records ---
# Row 20 becomes what row 22 used to be
$ws.Range("A20").Value = 131314800
$ws.Range("B20").Value = 57884
$ws.Range("E20").Value = 100109
$ws.Range("F20").Value = "Tretåig hackspett"
$ws.Range("G20").Value = "Picoides tridactylus"
$ws.Range("H20").Value = "(Linnaeus, 1758)"
$ws.Range("J20").Value = ""
$ws.Range("M20").Value = "äldre spår"
$ws.Range("Q20").Value = 497562
$ws.Range("R20").Value = 6980376
$ws.Range("AC20").Value = "Ringhack, äldre, på gran."
$ws.Range("AF20").Value = ""

# Row 22 becomes what row 20 used to be
$ws.Range("A22").Value = 131314710
$ws.Range("B22").Value = 79245
$ws.Range("E22").Value = 6425
$ws.Range("F22").Value = "Garnlav"
$ws.Range("G22").Value = "Alectoria sarmentosa"
$ws.Range("H22").Value = "(Ach.) Ach."
$ws.Range("M22").Value = ""
$ws.Range("Q22").Value = 497642
$ws.Range("R22").Value = 6980349
$ws.Range("AC22").Value = "Ringhack, äldre, på gran."
